# Desarrollo de proceso 100
#
# Refresh the "Precio" values that were re-scraped for this batch. The
# price cells are stored as plain text (not numbers), so each update is
# written via a temporary text-returning formula and then flattened back
# to a static value (Copy / PasteSpecial values-only) - this keeps the
# cell's type as text instead of letting a bare numeric-looking string be
# auto-coerced into a Number by Excel's input parser.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# "Productos" sheet - main product list
$wsProductos = $wb.Worksheets.Item("Productos")

$productosUpdates = [ordered]@{
    "C2"  = "7471"
    "C3"  = "27424"
    "C4"  = "11211"
    "C5"  = "84143"
    "C6"  = "186997"
    "C8"  = "14456"
    "C9"  = "14951"
    "C10" = "17195"
    "C12" = "466580"
    "C15" = "13455"
    "C16" = "13362"
    "C17" = "39262"
    "C18" = "71025"
    "C20" = "29640"
}

foreach ($addr in $productosUpdates.Keys) {
    Set-TextValue $wsProductos.Range($addr) $productosUpdates[$addr]
}

# "Productos Filtrados" sheet - filtered subset, same price figures
$wsFiltrados = $wb.Worksheets.Item("Productos Filtrados")

$filtradosUpdates = [ordered]@{
    "B2"  = "27424"
    "B3"  = "11211"
    "B4"  = "14456"
    "B5"  = "14951"
    "B6"  = "17195"
    "B7"  = "13455"
    "B8"  = "13362"
    "B9"  = "39262"
    "B10" = "29640"
}

foreach ($addr in $filtradosUpdates.Keys) {
    Set-TextValue $wsFiltrados.Range($addr) $filtradosUpdates[$addr]
}

$excel.CutCopyMode = $false
